$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1755.25
$ws.Range("I58").Value = 1755.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 5265.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -5115.75
$ws.Range("H64").Value = 250004000
$ws.Range("I64").Value = 250004000
$ws.Range("K64").Value = 250004000
$ws.Range("M64").Value = -250003752
$ws.Range("H67").Value = 250004000
$ws.Range("I67").Value = 250004000
$ws.Range("K67").Value = 250004000
$ws.Range("M67").Value = -250003142
$ws.Range("H98").Value = 3832.0356
$ws.Range("I98").Value = 3396.2307
$ws.Range("K98").Value = 3396.2307
$ws.Range("M98").Value = -1898.2307
$ws.Range("H122").Value = 3832.0356
$ws.Range("I122").Value = 3396.2307
$ws.Range("K122").Value = 10188.6921
$ws.Range("M122").Value = -7738.6921
$ws.Range("H137").Value = 9033.179
$ws.Range("I137").Value = 1267.762
$ws.Range("J137").Value = 32329.428
$ws.Range("K137").Value = 3803.286
$ws.Range("L137").Value = 96988.284
$ws.Range("M137").Value = -1253.286
$ws.Range("N137").Value = -102088.284
$ws.Range("H141").Value = 3038.08
$ws.Range("I141").Value = 3298.5557
$ws.Range("K141").Value = 9895.667099999999
$ws.Range("M141").Value = -4715.667099999999
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2807.5
$ws.Range("I2").Value = 2449
$ws.Range("J2").Value = 3524.5
$ws.Range("K2").Value = 2449
$ws.Range("L2").Value = 3524.5
$ws.Range("M2").Value = -2336
$ws.Range("N2").Value = -3750.5
$ws.Range("H116").Value = 2807.5
$ws.Range("I116").Value = 2449
$ws.Range("J116").Value = 3524.5
$ws.Range("K116").Value = 2449
$ws.Range("L116").Value = 3524.5
$ws.Range("M116").Value = -155
$ws.Range("N116").Value = -8112.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2807.5
$ws.Range("I3").Value = 2449
$ws.Range("J3").Value = 3524.5
$ws.Range("K3").Value = 2449
$ws.Range("L3").Value = 3524.5
$ws.Range("M3").Value = -2335
$ws.Range("N3").Value = -3752.5
$ws.Range("H20").Value = 5863193
$ws.Range("I20").Value = 9807689
$ws.Range("K20").Value = 9807689
$ws.Range("M20").Value = -9807442
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9203.073
$ws.Range("I31").Value = 1365.2916
$ws.Range("K31").Value = 1365.2916
$ws.Range("M31").Value = -1070.2916
$ws.Range("H34").Value = 9203.073
$ws.Range("I34").Value = 1365.2916
$ws.Range("K34").Value = 1365.2916
$ws.Range("M34").Value = -1163.2916
$ws.Range("H58").Value = 10711.173
$ws.Range("I58").Value = 3882.639
$ws.Range("J58").Value = 26075.375
$ws.Range("K58").Value = 3882.639
$ws.Range("L58").Value = 26075.375
$ws.Range("M58").Value = -3679.639
$ws.Range("N58").Value = -26481.375
$ws.Range("H76").Value = 7555.4443
$ws.Range("I76").Value = 7555.4443
$ws.Range("K76").Value = 7555.4443
$ws.Range("M76").Value = -7240.4443
$ws.Range("H79").Value = 7555.4443
$ws.Range("I79").Value = 7555.4443
$ws.Range("K79").Value = 7555.4443
$ws.Range("M79").Value = -6463.4443
$ws.Range("H86").Value = 9694.591
$ws.Range("I86").Value = 11247.5625
$ws.Range("K86").Value = 11247.5625
$ws.Range("M86").Value = -10124.5625
$ws.Range("H89").Value = 9694.591
$ws.Range("I89").Value = 11247.5625
$ws.Range("K89").Value = 56237.8125
$ws.Range("M89").Value = -50621.8125
$ws.Range("H99").Value = 3987100.2
$ws.Range("J99").Value = 4010749.5
$ws.Range("L99").Value = 4010749.5
$ws.Range("N99").Value = -4013745.5
$ws.Range("H126").Value = 3987100.2
$ws.Range("J126").Value = 4010749.5
$ws.Range("L126").Value = 12032248.5
$ws.Range("N126").Value = -12037188.5
$ws.Range("H132").Value = 2978.7778
$ws.Range("I132").Value = 2015.2
$ws.Range("K132").Value = 6045.6
$ws.Range("M132").Value = -3515.6
$ws.Range("H136").Value = 10711.173
$ws.Range("I136").Value = 3882.639
$ws.Range("J136").Value = 26075.375
$ws.Range("K136").Value = 11647.917
$ws.Range("L136").Value = 78226.125
$ws.Range("M136").Value = -9097.917000000001
$ws.Range("N136").Value = -83326.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 145.25
$ws.Range("I33").Value = 201
$ws.Range("J33").Value = 126.666664
$ws.Range("K33").Value = 1206
$ws.Range("L33").Value = 759.999984
$ws.Range("M33").Value = -923
$ws.Range("N33").Value = -1325.999984
$ws.Range("H38").Value = 98.25
$ws.Range("I38").Value = 44.142857
$ws.Range("J38").Value = 174
$ws.Range("K38").Value = 132.428571
$ws.Range("L38").Value = 522
$ws.Range("M38").Value = 214.571429
$ws.Range("N38").Value = -1216
$ws.Range("H44").Value = 1458.6666
$ws.Range("I44").Value = 1280
$ws.Range("J44").Value = 2352
$ws.Range("K44").Value = 3840
$ws.Range("L44").Value = 7056
$ws.Range("M44").Value = -3442
$ws.Range("N44").Value = -7852
$ws.Range("H47").Value = 14844019
$ws.Range("H48").Value = 4978265
$ws.Range("I48").Value = 2489133
$ws.Range("J48").Value = 7467397
$ws.Range("K48").Value = 7467399
$ws.Range("L48").Value = 22402191
$ws.Range("M48").Value = -7467149
$ws.Range("N48").Value = -22402691
$ws.Range("H54").Value = 5821.2
$ws.Range("I54").Value = 21212
$ws.Range("J54").Value = 4111.1113
$ws.Range("K54").Value = 63636
$ws.Range("L54").Value = 12333.3339
$ws.Range("M54").Value = -63077
$ws.Range("N54").Value = -13451.3339
$ws.Range("H60").Value = 2802.25
$ws.Range("I60").Value = 3403
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 10209
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -9958
$ws.Range("N60").Value = -3502
$ws.Range("H61").Value = 373
$ws.Range("I61").Value = 373
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1119
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -904
$ws.Range("H63").Value = 10000
$ws.Range("J63").Value = 10000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31498
$ws.Range("H66").Value = 10000
$ws.Range("J66").Value = 10000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -97488
$ws.Range("H122").Value = 7688467
$ws.Range("J122").Value = 1494771.6
$ws.Range("L122").Value = 13452944.4
$ws.Range("N122").Value = -13457844.4
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4226.1875
$ws.Range("I61").Value = 3259.9167
$ws.Range("J61").Value = 7125
$ws.Range("K61").Value = 3259.9167
$ws.Range("L61").Value = 7125
$ws.Range("M61").Value = -3057.9167
$ws.Range("N61").Value = -7529
$ws.Range("H68").Value = 1599.5
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 2000
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 1599.5
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 10000
$ws.Range("N71").Value = -17488
$ws.Range("H93").Value = 41671280
$ws.Range("I93").Value = 66670908
$ws.Range("J93").Value = 5233.778
$ws.Range("K93").Value = 66670908
$ws.Range("L93").Value = 5233.778
$ws.Range("M93").Value = -66669660
$ws.Range("N93").Value = -7729.778
$ws.Range("H100").Value = 4520.1763
$ws.Range("J100").Value = 2998.4
$ws.Range("L100").Value = 2998.4
$ws.Range("N100").Value = -4080.4
$ws.Range("H113").Value = 4226.1875
$ws.Range("I113").Value = 3259.9167
$ws.Range("J113").Value = 7125
$ws.Range("K113").Value = 3259.9167
$ws.Range("L113").Value = 7125
$ws.Range("M113").Value = -1089.9167
$ws.Range("N113").Value = -11465
$ws.Range("H132").Value = 1679148
$ws.Range("I132").Value = 2033.2222
$ws.Range("K132").Value = 6099.6666
$ws.Range("M132").Value = -3569.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 49833
$ws.Range("J74").Value = 49833
$ws.Range("L74").Value = 49833
$ws.Range("N74").Value = -51705
$ws.Range("H77").Value = 49833
$ws.Range("J77").Value = 49833
$ws.Range("L77").Value = 149499
$ws.Range("N77").Value = -158859
$ws.Range("H132").Value = 9616.888999999999
$ws.Range("I132").Value = 3208.476
$ws.Range("K132").Value = 9625.428
$ws.Range("M132").Value = -7095.428
